$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Phase 1: hold ws2 (总计) original formats in an unused scratch area of ws2 itself ---
$ws2.Range("B1:D1").Copy()
$ws2.Range("B100:D100").PasteSpecial(-4122)
$ws2.Range("A2").Copy()
$ws2.Range("A101").PasteSpecial(-4122)

# --- Phase 2: transplant ws1 (2020-Q4) formats onto ws2 final layout ---
$ws1.Range("B1:H1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2:A42").Copy()
$ws2.Range("A2:A42").PasteSpecial(-4122)

# --- Phase 3: transplant held ws2 (总计) formats onto ws1 final layout ---
$ws2.Range("B100:D100").Copy()
$ws1.Range("B1:D1").PasteSpecial(-4122)
$ws2.Range("A101").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

# --- Phase 4: clear scratch holding area + clear leftover ws1 cells outside new A1:D2 layout ---
$ws2.Range("A100:D101").Clear()
$ws1.Range("E1:H42").Clear()
$ws1.Range("A3:D42").Clear()

# --- Phase 5: rename sheets ---
$ws1.Name = "总计"
$ws2.Name = "2020-Q4"

# --- Phase 6: write values ---
# 总计 (ws1) values
$ws1.Cells.Item(1,2).Value = "日期"
$ws1.Cells.Item(1,3).Value = "持有数量(只)"
$ws1.Cells.Item(1,4).Value = "持有市值(亿元)"
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2020-Q4"
$ws1.Cells.Item(2,3).Value = 41
$ws1.Cells.Item(2,4).Value = 25.69

# 2020-Q4 (ws2) values
$ws2.Cells.Item(1,2).Value = "基金代码"
$ws2.Cells.Item(1,3).Value = "基金名称"
$ws2.Cells.Item(1,4).Value = "基金金额"
$ws2.Cells.Item(1,5).Value = "股票总仓位"
$ws2.Cells.Item(1,6).Value = "仓位占比"
$ws2.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1,8).Value = "仓位排名"
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "'001875"
$ws2.Cells.Item(2,3).Value = "前海开源沪港深优势精选灵活配置混合"
$ws2.Cells.Item(2,4).Value = "'73.55"
$ws2.Cells.Item(2,5).Value = "'94.53"
$ws2.Cells.Item(2,6).Value = "'8.64"
$ws2.Cells.Item(2,7).Value = "'6.3547"
$ws2.Cells.Item(2,8).Value = 9
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "'005379"
$ws2.Cells.Item(3,3).Value = "汇添富价值创造定期开放混合"
$ws2.Cells.Item(3,4).Value = "'95.71"
$ws2.Cells.Item(3,5).Value = "'93.98"
$ws2.Cells.Item(3,6).Value = "'4.98"
$ws2.Cells.Item(3,7).Value = "'4.7664"
$ws2.Cells.Item(3,8).Value = 6
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).Value = "'001837"
$ws2.Cells.Item(4,3).Value = "前海开源沪港深蓝筹精选灵活配置混合"
$ws2.Cells.Item(4,4).Value = "'37.41"
$ws2.Cells.Item(4,5).Value = "'94.05"
$ws2.Cells.Item(4,6).Value = "'9.05"
$ws2.Cells.Item(4,7).Value = "'3.3856"
$ws2.Cells.Item(4,8).Value = 6
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).Value = "'501065"
$ws2.Cells.Item(5,3).Value = "汇添富经典成长定期开放混合"
$ws2.Cells.Item(5,4).Value = "'14.25"
$ws2.Cells.Item(5,5).Value = "'73.59"
$ws2.Cells.Item(5,6).Value = "'9.61"
$ws2.Cells.Item(5,7).Value = "'1.3694"
$ws2.Cells.Item(5,8).Value = 3
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).Value = "'006408"
$ws2.Cells.Item(6,3).Value = "汇添富消费升级混合"
$ws2.Cells.Item(6,4).Value = "'37.46"
$ws2.Cells.Item(6,5).Value = "'94.25"
$ws2.Cells.Item(6,6).Value = "'3.54"
$ws2.Cells.Item(6,7).Value = "'1.3261"
$ws2.Cells.Item(6,8).Value = 10
$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,2).Value = "'009007"
$ws2.Cells.Item(7,3).Value = "兴全沪港深两年持有期混合"
$ws2.Cells.Item(7,4).Value = "'32.05"
$ws2.Cells.Item(7,5).Value = "'90.67"
$ws2.Cells.Item(7,6).Value = "'3.89"
$ws2.Cells.Item(7,7).Value = "'1.2467"
$ws2.Cells.Item(7,8).Value = 2
$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,2).Value = "'005644"
$ws2.Cells.Item(8,3).Value = "广发沪港深行业龙头混合"
$ws2.Cells.Item(8,4).Value = "'13.29"
$ws2.Cells.Item(8,5).Value = "'94.52"
$ws2.Cells.Item(8,6).Value = "'8.71"
$ws2.Cells.Item(8,7).Value = "'1.1576"
$ws2.Cells.Item(8,8).Value = 3
$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,2).Value = "'009931"
$ws2.Cells.Item(9,3).Value = "淳厚欣享一年持有期混合A"
$ws2.Cells.Item(9,4).Value = "'16.79"
$ws2.Cells.Item(9,5).Value = "'92.51"
$ws2.Cells.Item(9,6).Value = "'4.98"
$ws2.Cells.Item(9,7).Value = "'0.8361"
$ws2.Cells.Item(9,8).Value = 2
$ws2.Cells.Item(10,1).Value = 8
$ws2.Cells.Item(10,2).Value = "'008025"
$ws2.Cells.Item(10,3).Value = "汇添富稳健增长混合A"
$ws2.Cells.Item(10,4).Value = "'20.70"
$ws2.Cells.Item(10,5).Value = "'39.49"
$ws2.Cells.Item(10,6).Value = "'3.20"
$ws2.Cells.Item(10,7).Value = "'0.6624"
$ws2.Cells.Item(10,8).Value = 5
$ws2.Cells.Item(11,1).Value = 9
$ws2.Cells.Item(11,2).Value = "'001874"
$ws2.Cells.Item(11,3).Value = "前海开源沪港深价值精选灵活配置混合"
$ws2.Cells.Item(11,4).Value = "'6.87"
$ws2.Cells.Item(11,5).Value = "'94.88"
$ws2.Cells.Item(11,6).Value = "'8.94"
$ws2.Cells.Item(11,7).Value = "'0.6142"
$ws2.Cells.Item(11,8).Value = 8
$ws2.Cells.Item(12,1).Value = 10
$ws2.Cells.Item(12,2).Value = "'002387"
$ws2.Cells.Item(12,3).Value = "工银瑞信沪港深股票A"
$ws2.Cells.Item(12,4).Value = "'13.26"
$ws2.Cells.Item(12,5).Value = "'93.04"
$ws2.Cells.Item(12,6).Value = "'4.51"
$ws2.Cells.Item(12,7).Value = "'0.5980"
$ws2.Cells.Item(12,8).Value = 4
$ws2.Cells.Item(13,1).Value = 11
$ws2.Cells.Item(13,2).Value = "'008186"
$ws2.Cells.Item(13,3).Value = "淳厚信睿核心精选混合A"
$ws2.Cells.Item(13,4).Value = "'8.68"
$ws2.Cells.Item(13,5).Value = "'93.67"
$ws2.Cells.Item(13,6).Value = "'5.19"
$ws2.Cells.Item(13,7).Value = "'0.4505"
$ws2.Cells.Item(13,8).Value = 2
$ws2.Cells.Item(14,1).Value = 12
$ws2.Cells.Item(14,2).Value = "'005583"
$ws2.Cells.Item(14,3).Value = "易方达港股通红利灵活配置混合"
$ws2.Cells.Item(14,4).Value = "'7.87"
$ws2.Cells.Item(14,5).Value = "'91.52"
$ws2.Cells.Item(14,6).Value = "'5.01"
$ws2.Cells.Item(14,7).Value = "'0.3943"
$ws2.Cells.Item(14,8).Value = 7
$ws2.Cells.Item(15,1).Value = 13
$ws2.Cells.Item(15,2).Value = "'002653"
$ws2.Cells.Item(15,3).Value = "泰康沪港深精选灵活配置混合"
$ws2.Cells.Item(15,4).Value = "'10.79"
$ws2.Cells.Item(15,5).Value = "'93.04"
$ws2.Cells.Item(15,6).Value = "'3.31"
$ws2.Cells.Item(15,7).Value = "'0.3571"
$ws2.Cells.Item(15,8).Value = 5
$ws2.Cells.Item(16,1).Value = 14
$ws2.Cells.Item(16,2).Value = "'002332"
$ws2.Cells.Item(16,3).Value = "汇丰晋信沪港深股票A"
$ws2.Cells.Item(16,4).Value = "'14.55"
$ws2.Cells.Item(16,5).Value = "'94.49"
$ws2.Cells.Item(16,6).Value = "'2.41"
$ws2.Cells.Item(16,7).Value = "'0.3507"
$ws2.Cells.Item(16,8).Value = 10
$ws2.Cells.Item(17,1).Value = 15
$ws2.Cells.Item(17,2).Value = "'501021"
$ws2.Cells.Item(17,3).Value = "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)A"
$ws2.Cells.Item(17,4).Value = "'12.34"
$ws2.Cells.Item(17,5).Value = "'94.72"
$ws2.Cells.Item(17,6).Value = "'2.05"
$ws2.Cells.Item(17,7).Value = "'0.2530"
$ws2.Cells.Item(17,8).Value = 8
$ws2.Cells.Item(18,1).Value = 16
$ws2.Cells.Item(18,2).Value = "'009939"
$ws2.Cells.Item(18,3).Value = "淳厚欣享一年持有期混合C"
$ws2.Cells.Item(18,4).Value = "'4.33"
$ws2.Cells.Item(18,5).Value = "'92.51"
$ws2.Cells.Item(18,6).Value = "'4.98"
$ws2.Cells.Item(18,7).Value = "'0.2156"
$ws2.Cells.Item(18,8).Value = 2
$ws2.Cells.Item(19,1).Value = 17
$ws2.Cells.Item(19,2).Value = "'008026"
$ws2.Cells.Item(19,3).Value = "汇添富稳健增长混合C"
$ws2.Cells.Item(19,4).Value = "'6.16"
$ws2.Cells.Item(19,5).Value = "'39.49"
$ws2.Cells.Item(19,6).Value = "'3.20"
$ws2.Cells.Item(19,7).Value = "'0.1971"
$ws2.Cells.Item(19,8).Value = 5
$ws2.Cells.Item(20,1).Value = 18
$ws2.Cells.Item(20,2).Value = "'001901"
$ws2.Cells.Item(20,3).Value = "前海开源沪港深隆鑫灵活配置混合A"
$ws2.Cells.Item(20,4).Value = "'5.00"
$ws2.Cells.Item(20,5).Value = "'47.42"
$ws2.Cells.Item(20,6).Value = "'3.45"
$ws2.Cells.Item(20,7).Value = "'0.1725"
$ws2.Cells.Item(20,8).Value = 7
$ws2.Cells.Item(21,1).Value = 19
$ws2.Cells.Item(21,2).Value = "'008187"
$ws2.Cells.Item(21,3).Value = "淳厚信睿核心精选混合C"
$ws2.Cells.Item(21,4).Value = "'2.78"
$ws2.Cells.Item(21,5).Value = "'93.67"
$ws2.Cells.Item(21,6).Value = "'5.19"
$ws2.Cells.Item(21,7).Value = "'0.1443"
$ws2.Cells.Item(21,8).Value = 2
$ws2.Cells.Item(22,1).Value = 20
$ws2.Cells.Item(22,2).Value = "'241001"
$ws2.Cells.Item(22,3).Value = "华宝海外中国混合(QDII)"
$ws2.Cells.Item(22,4).Value = "'1.03"
$ws2.Cells.Item(22,5).Value = "'94.45"
$ws2.Cells.Item(22,6).Value = "'8.03"
$ws2.Cells.Item(22,7).Value = "'0.0827"
$ws2.Cells.Item(22,8).Value = 2
$ws2.Cells.Item(23,1).Value = 21
$ws2.Cells.Item(23,2).Value = "'007291"
$ws2.Cells.Item(23,3).Value = "汇丰晋信港股通双核策略混合"
$ws2.Cells.Item(23,4).Value = "'2.91"
$ws2.Cells.Item(23,5).Value = "'94.38"
$ws2.Cells.Item(23,6).Value = "'2.70"
$ws2.Cells.Item(23,7).Value = "'0.0786"
$ws2.Cells.Item(23,8).Value = 9
$ws2.Cells.Item(24,1).Value = 22
$ws2.Cells.Item(24,2).Value = "'002443"
$ws2.Cells.Item(24,3).Value = "前海开源沪港深龙头精选灵活配置混合"
$ws2.Cells.Item(24,4).Value = "'0.97"
$ws2.Cells.Item(24,5).Value = "'94.89"
$ws2.Cells.Item(24,6).Value = "'7.91"
$ws2.Cells.Item(24,7).Value = "'0.0767"
$ws2.Cells.Item(24,8).Value = 5
$ws2.Cells.Item(25,1).Value = 23
$ws2.Cells.Item(25,2).Value = "'008381"
$ws2.Cells.Item(25,3).Value = "前海开源新兴产业混合"
$ws2.Cells.Item(25,4).Value = "'0.88"
$ws2.Cells.Item(25,5).Value = "'94.90"
$ws2.Cells.Item(25,6).Value = "'8.72"
$ws2.Cells.Item(25,7).Value = "'0.0767"
$ws2.Cells.Item(25,8).Value = 4
$ws2.Cells.Item(26,1).Value = 24
$ws2.Cells.Item(26,2).Value = "'003580"
$ws2.Cells.Item(26,3).Value = "泰康沪港深价值优选灵活配置混合"
$ws2.Cells.Item(26,4).Value = "'1.54"
$ws2.Cells.Item(26,5).Value = "'93.73"
$ws2.Cells.Item(26,6).Value = "'4.29"
$ws2.Cells.Item(26,7).Value = "'0.0661"
$ws2.Cells.Item(26,8).Value = 4
$ws2.Cells.Item(27,1).Value = 25
$ws2.Cells.Item(27,2).Value = "'007151"
$ws2.Cells.Item(27,3).Value = "前海开源沪港深聚瑞混合"
$ws2.Cells.Item(27,4).Value = "'0.68"
$ws2.Cells.Item(27,5).Value = "'93.57"
$ws2.Cells.Item(27,6).Value = "'9.72"
$ws2.Cells.Item(27,7).Value = "'0.0661"
$ws2.Cells.Item(27,8).Value = 2
$ws2.Cells.Item(28,1).Value = 26
$ws2.Cells.Item(28,2).Value = "'006049"
$ws2.Cells.Item(28,3).Value = "恒越研究精选混合A/B"
$ws2.Cells.Item(28,4).Value = "'2.11"
$ws2.Cells.Item(28,5).Value = "'82.80"
$ws2.Cells.Item(28,6).Value = "'2.97"
$ws2.Cells.Item(28,7).Value = "'0.0627"
$ws2.Cells.Item(28,8).Value = 9
$ws2.Cells.Item(29,1).Value = 27
$ws2.Cells.Item(29,2).Value = "'007192"
$ws2.Cells.Item(29,3).Value = "恒越研究精选混合C"
$ws2.Cells.Item(29,4).Value = "'2.11"
$ws2.Cells.Item(29,5).Value = "'82.80"
$ws2.Cells.Item(29,6).Value = "'2.97"
$ws2.Cells.Item(29,7).Value = "'0.0627"
$ws2.Cells.Item(29,8).Value = 9
$ws2.Cells.Item(30,1).Value = 28
$ws2.Cells.Item(30,2).Value = "'519601"
$ws2.Cells.Item(30,3).Value = "海富通中国海外精选混合(QDII)"
$ws2.Cells.Item(30,4).Value = "'1.32"
$ws2.Cells.Item(30,5).Value = "'93.40"
$ws2.Cells.Item(30,6).Value = "'4.17"
$ws2.Cells.Item(30,7).Value = "'0.0550"
$ws2.Cells.Item(30,8).Value = 7
$ws2.Cells.Item(31,1).Value = 29
$ws2.Cells.Item(31,2).Value = "'005534"
$ws2.Cells.Item(31,3).Value = "华夏新时代灵活配置混合（QDII）"
$ws2.Cells.Item(31,4).Value = "'1.37"
$ws2.Cells.Item(31,5).Value = "'90.36"
$ws2.Cells.Item(31,6).Value = "'2.70"
$ws2.Cells.Item(31,7).Value = "'0.0370"
$ws2.Cells.Item(31,8).Value = 8
$ws2.Cells.Item(32,1).Value = 30
$ws2.Cells.Item(32,2).Value = "'007280"
$ws2.Cells.Item(32,3).Value = "上投摩根日本精选股票（QDII）"
$ws2.Cells.Item(32,4).Value = "'1.11"
$ws2.Cells.Item(32,5).Value = "'94.14"
$ws2.Cells.Item(32,6).Value = "'3.19"
$ws2.Cells.Item(32,7).Value = "'0.0354"
$ws2.Cells.Item(32,8).Value = 8
$ws2.Cells.Item(33,1).Value = 31
$ws2.Cells.Item(33,2).Value = "'002333"
$ws2.Cells.Item(33,3).Value = "汇丰晋信沪港深股票C"
$ws2.Cells.Item(33,4).Value = "'1.06"
$ws2.Cells.Item(33,5).Value = "'94.49"
$ws2.Cells.Item(33,6).Value = "'2.41"
$ws2.Cells.Item(33,7).Value = "'0.0255"
$ws2.Cells.Item(33,8).Value = 10
$ws2.Cells.Item(34,1).Value = 32
$ws2.Cells.Item(34,2).Value = "'006127"
$ws2.Cells.Item(34,3).Value = "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)C"
$ws2.Cells.Item(34,4).Value = "'1.14"
$ws2.Cells.Item(34,5).Value = "'94.72"
$ws2.Cells.Item(34,6).Value = "'2.05"
$ws2.Cells.Item(34,7).Value = "'0.0234"
$ws2.Cells.Item(34,8).Value = 8
$ws2.Cells.Item(35,1).Value = 33
$ws2.Cells.Item(35,2).Value = "'007132"
$ws2.Cells.Item(35,3).Value = "长城港股通价值精选多策略混合"
$ws2.Cells.Item(35,4).Value = "'0.54"
$ws2.Cells.Item(35,5).Value = "'88.86"
$ws2.Cells.Item(35,6).Value = "'4.11"
$ws2.Cells.Item(35,7).Value = "'0.0222"
$ws2.Cells.Item(35,8).Value = 6
$ws2.Cells.Item(36,1).Value = 34
$ws2.Cells.Item(36,2).Value = "'005255"
$ws2.Cells.Item(36,3).Value = "浦银安盛港股通量化优选灵活配置混合"
$ws2.Cells.Item(36,4).Value = "'0.41"
$ws2.Cells.Item(36,5).Value = "'92.69"
$ws2.Cells.Item(36,6).Value = "'3.42"
$ws2.Cells.Item(36,7).Value = "'0.0140"
$ws2.Cells.Item(36,8).Value = 7
$ws2.Cells.Item(37,1).Value = 35
$ws2.Cells.Item(37,2).Value = "'080006"
$ws2.Cells.Item(37,3).Value = "长盛环球行业混合(QDII)"
$ws2.Cells.Item(37,4).Value = "'0.25"
$ws2.Cells.Item(37,5).Value = "'91.90"
$ws2.Cells.Item(37,6).Value = "'4.77"
$ws2.Cells.Item(37,7).Value = "'0.0119"
$ws2.Cells.Item(37,8).Value = 5
$ws2.Cells.Item(38,1).Value = 36
$ws2.Cells.Item(38,2).Value = "'006816"
$ws2.Cells.Item(38,3).Value = "泰康中证港股通地产指数A"
$ws2.Cells.Item(38,4).Value = "'0.25"
$ws2.Cells.Item(38,5).Value = "'92.44"
$ws2.Cells.Item(38,6).Value = "'4.22"
$ws2.Cells.Item(38,7).Value = "'0.0106"
$ws2.Cells.Item(38,8).Value = 8
$ws2.Cells.Item(39,1).Value = 37
$ws2.Cells.Item(39,2).Value = "'007512"
$ws2.Cells.Item(39,3).Value = "工银瑞信沪港深股票C"
$ws2.Cells.Item(39,4).Value = "'0.23"
$ws2.Cells.Item(39,5).Value = "'93.04"
$ws2.Cells.Item(39,6).Value = "'4.51"
$ws2.Cells.Item(39,7).Value = "'0.0104"
$ws2.Cells.Item(39,8).Value = 4
$ws2.Cells.Item(40,1).Value = 38
$ws2.Cells.Item(40,2).Value = "'519602"
$ws2.Cells.Item(40,3).Value = "海富通大中华精选混合QDII"
$ws2.Cells.Item(40,4).Value = "'0.24"
$ws2.Cells.Item(40,5).Value = "'93.35"
$ws2.Cells.Item(40,6).Value = "'4.29"
$ws2.Cells.Item(40,7).Value = "'0.0103"
$ws2.Cells.Item(40,8).Value = 7
$ws2.Cells.Item(41,1).Value = 39
$ws2.Cells.Item(41,2).Value = "'006817"
$ws2.Cells.Item(41,3).Value = "泰康中证港股通地产指数C"
$ws2.Cells.Item(41,4).Value = "'0.19"
$ws2.Cells.Item(41,5).Value = "'92.44"
$ws2.Cells.Item(41,6).Value = "'4.22"
$ws2.Cells.Item(41,7).Value = "'0.0080"
$ws2.Cells.Item(41,8).Value = 8
$ws2.Cells.Item(42,1).Value = 40
$ws2.Cells.Item(42,2).Value = "'001902"
$ws2.Cells.Item(42,3).Value = "前海开源沪港深隆鑫灵活配置混合C"
$ws2.Cells.Item(42,4).Value = "'0.04"
$ws2.Cells.Item(42,5).Value = "'47.42"
$ws2.Cells.Item(42,6).Value = "'3.45"
$ws2.Cells.Item(42,7).Value = "'0.0014"
$ws2.Cells.Item(42,8).Value = 7

# --- Phase 7: page margins (points = inches * 72) ---
$ws1.PageSetup.LeftMargin = 54
$ws1.PageSetup.RightMargin = 54
$ws1.PageSetup.TopMargin = 72
$ws1.PageSetup.BottomMargin = 72
$ws1.PageSetup.HeaderMargin = 36
$ws1.PageSetup.FooterMargin = 36

$ws2.PageSetup.LeftMargin = 50.4
$ws2.PageSetup.RightMargin = 50.4
$ws2.PageSetup.TopMargin = 54
$ws2.PageSetup.BottomMargin = 54
$ws2.PageSetup.HeaderMargin = 21.6
$ws2.PageSetup.FooterMargin = 21.6

# --- Phase 8: keep "2020-Q4" as the selected/active tab ---
$ws2.Activate()
$ws1.Range("A1").Select()
$ws2.Range("A1").Select()